# Lab01_ReviewReport.xlsx - "Fix some things for the report"
#
# - Architect. Design Phase Defects (sheet 2): the "Doc. page/line" column
#   (D10:D12) was a placeholder numeric 1 - replace with the real value
#   "Diagram".
# - Coding Phase Defects (sheet 3): the "Doc. page/line" entries get more
#   precise line references.
# - Update the last-used cell selection remembered on sheet 2 and sheet 3.

$wb = $excel.ActiveWorkbook

# --- Window / view state ---------------------------------------------------
# Best-effort: mirror the saved window geometry (yWindow 1950 -> 2400) and
# the scrolled-tab position. Not every host surfaces these as writable
# properties, so failures here are non-fatal to the rest of the edit.
try {
    $win = $excel.ActiveWindow
    $win.Top = 2400
    $win.Left = 0
    $win.Width = 14160
    $win.Height = 8265
    $win.ScrollWorkbookTabs(1, 1) | Out-Null
} catch {
}

# --- Architect. Design Phase Defects -------------------------------------
$wsArch = $wb.Worksheets.Item("Architect. Design Phase Defects")

$wsArch.Range("D10").Value = "Diagram"
$wsArch.Range("D11").Value = "Diagram"
$wsArch.Range("D12").Value = "Diagram"

# remembered selection for this sheet moves to D12
$wsArch.Range("D12").Select()

# --- Coding Phase Defects --------------------------------------------------
$wsCode = $wb.Worksheets.Item("Coding Phase Defects")

$wsCode.Range("D10").Value = "File Repos, lines 18-19"
$wsCode.Range("D11").Value = "Activity, line 12"

# remembered selection for this sheet moves to H12, and this sheet stays
# the active tab (it was already the active tab before the edit)
$wsCode.Range("H12").Select()
